$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = @(4, 0, 1)
    3  = @(1, 2, 6)
    4  = @(10, 5, 2)
    5  = @(7, 5, 7)
    6  = @(2, 6, 2)
    7  = @(6, 2, 6)
    8  = @(4, 0, 1)
    9  = @(6, 9, 6)
    10 = @(10, 4, 9)
    11 = @(3, 10, 7)
    12 = @(10, 1, 5)
    13 = @(5, 10, 7)
    14 = @(10, 7, 2)
    15 = @(3, 8, 8)
    16 = @(5, 6, 1)
    17 = @(8, 1, 7)
    18 = @(10, 9, 5)
    19 = @(2, 4, 8)
    20 = @(6, 0, 2)
    21 = @(4, 1, 10)
}

foreach ($row in $values.Keys) {
    $triple = $values[$row]
    $ws.Range("D" + $row).Value = $triple[0]
    $ws.Range("E" + $row).Value = $triple[1]
    $ws.Range("F" + $row).Value = $triple[2]
}
